# Scheduled-runner market data refresh: update computed price/profit
# columns (H-N) for affected leve rows across the ALC/ARM/CRP/CUL/GSM/
# LTW/WVR sheets, per the upstream Sheets diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 140.3
$ws.Range("I6").Value = 56.625
$ws.Range("K6").Value = 169.875
$ws.Range("M6").Value = -57.875

$ws.Range("H28").Value = 540.4074000000001
$ws.Range("I28").Value = 565.56525
$ws.Range("J28").Value = 395.75
$ws.Range("K28").Value = 565.56525
$ws.Range("L28").Value = 395.75
$ws.Range("M28").Value = -80.56524999999999
$ws.Range("N28").Value = -1365.75

$ws.Range("H107").Value = 627.2727
$ws.Range("I107").Value = 670.2778
$ws.Range("K107").Value = 670.2778
$ws.Range("M107").Value = 1249.7222

$ws.Range("H116").Value = 103699.57
$ws.Range("I116").Value = 178223.75
$ws.Range("J116").Value = 4334
$ws.Range("K116").Value = 178223.75
$ws.Range("L116").Value = 4334
$ws.Range("M116").Value = -174781.75
$ws.Range("N116").Value = -11218

$ws.Range("H137").Value = 29659.838
$ws.Range("I137").Value = 44642.332
$ws.Range("J137").Value = 1999.8462
$ws.Range("K137").Value = 133926.996
$ws.Range("L137").Value = 5999.5386
$ws.Range("M137").Value = -131376.996
$ws.Range("N137").Value = -11099.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1146.3055
$ws.Range("I110").Value = 1146.88
$ws.Range("K110").Value = 1146.88
$ws.Range("M110").Value = 898.1199999999999

$ws.Range("H132").Value = 3010.7285
$ws.Range("I132").Value = 2704.617
$ws.Range("J132").Value = 3636.261
$ws.Range("K132").Value = 8113.851000000001
$ws.Range("L132").Value = 10908.783
$ws.Range("M132").Value = -5583.851000000001
$ws.Range("N132").Value = -15968.783

$ws.Range("H133").Value = 34000
$ws.Range("J133").Value = 34000
$ws.Range("L133").Value = 34000
$ws.Range("N133").Value = -39060

$ws.Range("H134").Value = 39326.715
$ws.Range("J134").Value = 39326.715
$ws.Range("L134").Value = 39326.715
$ws.Range("N134").Value = -49466.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2847.218
$ws.Range("I31").Value = 1929.2368
$ws.Range("J31").Value = 3719.3
$ws.Range("K31").Value = 1929.2368
$ws.Range("L31").Value = 3719.3
$ws.Range("M31").Value = -1634.2368
$ws.Range("N31").Value = -4309.3

$ws.Range("H34").Value = 2847.218
$ws.Range("I34").Value = 1929.2368
$ws.Range("J34").Value = 3719.3
$ws.Range("K34").Value = 1929.2368
$ws.Range("L34").Value = 3719.3
$ws.Range("M34").Value = -1727.2368
$ws.Range("N34").Value = -4123.3

$ws.Range("H134").Value = 1906.9143
$ws.Range("I134").Value = 1118.4736
$ws.Range("K134").Value = 3355.4208
$ws.Range("M134").Value = -820.4207999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 604.96
$ws.Range("I5").Value = 384
$ws.Range("J5").Value = 1074.5
$ws.Range("K5").Value = 1152
$ws.Range("L5").Value = 3223.5
$ws.Range("M5").Value = -1040
$ws.Range("N5").Value = -3447.5

$ws.Range("H33").Value = 265
$ws.Range("I33").Value = 167.85715
$ws.Range("J33").Value = 378.33334
$ws.Range("K33").Value = 1007.1429
$ws.Range("L33").Value = 2270.00004
$ws.Range("M33").Value = -724.1428999999999
$ws.Range("N33").Value = -2836.00004

$ws.Range("H80").Value = 1378.1111
$ws.Range("J80").Value = 1437.875
$ws.Range("L80").Value = 4313.625
$ws.Range("N80").Value = -6185.625

$ws.Range("H83").Value = 1378.1111
$ws.Range("J83").Value = 1437.875
$ws.Range("L83").Value = 12940.875
$ws.Range("N83").Value = -22300.875

$ws.Range("H131").Value = 873.3692
$ws.Range("I131").Value = 493
$ws.Range("J131").Value = 912.05084
$ws.Range("K131").Value = 1479
$ws.Range("L131").Value = 2736.15252
$ws.Range("M131").Value = 3561
$ws.Range("N131").Value = -12816.15252

$ws.Range("H135").Value = 604.96
$ws.Range("I135").Value = 384
$ws.Range("J135").Value = 1074.5
$ws.Range("K135").Value = 3456
$ws.Range("L135").Value = 9670.5
$ws.Range("M135").Value = -921
$ws.Range("N135").Value = -14740.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1034.7778
$ws.Range("I97").Value = 679.875
$ws.Range("J97").Value = 1551
$ws.Range("K97").Value = 679.875
$ws.Range("L97").Value = 1551
$ws.Range("M97").Value = -183.875
$ws.Range("N97").Value = -2543

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 502.7
$ws.Range("I22").Value = 397.8
$ws.Range("J22").Value = 607.6
$ws.Range("K22").Value = 397.8
$ws.Range("L22").Value = 607.6
$ws.Range("M22").Value = -102.8
$ws.Range("N22").Value = -1197.6

$ws.Range("H27").Value = 502.7
$ws.Range("I27").Value = 397.8
$ws.Range("J27").Value = 607.6
$ws.Range("K27").Value = 397.8
$ws.Range("L27").Value = 607.6
$ws.Range("M27").Value = -290.8
$ws.Range("N27").Value = -821.6

$ws.Range("H43").Value = 5425
$ws.Range("I43").Value = 5100
$ws.Range("J43").Value = 5750
$ws.Range("K43").Value = 5100
$ws.Range("L43").Value = 5750
$ws.Range("M43").Value = -4907
$ws.Range("N43").Value = -6136

$ws.Range("H46").Value = 700.2
$ws.Range("I46").Value = 501
$ws.Range("J46").Value = 750
$ws.Range("K46").Value = 501
$ws.Range("L46").Value = 750
$ws.Range("M46").Value = -313
$ws.Range("N46").Value = -1126

$ws.Range("H68").Value = 62502044
$ws.Range("I68").Value = 111112730
$ws.Range("J68").Value = 2594.2856
$ws.Range("K68").Value = 111112730
$ws.Range("L68").Value = 2594.2856
$ws.Range("M68").Value = -111111981
$ws.Range("N68").Value = -4092.2856

$ws.Range("H71").Value = 62502044
$ws.Range("I71").Value = 111112730
$ws.Range("J71").Value = 2594.2856
$ws.Range("K71").Value = 555563650
$ws.Range("L71").Value = 12971.428
$ws.Range("M71").Value = -555559906
$ws.Range("N71").Value = -20459.428

$ws.Range("H82").Value = 1594.3182
$ws.Range("I82").Value = 1311.7273
$ws.Range("K82").Value = 1311.7273
$ws.Range("M82").Value = -950.7273

$ws.Range("H85").Value = 1594.3182
$ws.Range("I85").Value = 1311.7273
$ws.Range("K85").Value = 1311.7273
$ws.Range("M85").Value = -63.72730000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 4050
$ws.Range("I33").Value = 3200
$ws.Range("J33").Value = 4333.3335
$ws.Range("K33").Value = 3200
$ws.Range("L33").Value = 4333.3335
$ws.Range("M33").Value = -2950
$ws.Range("N33").Value = -4833.3335

$ws.Range("H36").Value = 4050
$ws.Range("I36").Value = 3200
$ws.Range("J36").Value = 4333.3335
$ws.Range("K36").Value = 3200
$ws.Range("L36").Value = 4333.3335
$ws.Range("M36").Value = -2950
$ws.Range("N36").Value = -4833.3335

$ws.Range("H40").Value = 8428.571
$ws.Range("J40").Value = 8428.571
$ws.Range("L40").Value = 8428.571
$ws.Range("N40").Value = -8726.571

$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -15756

$ws.Range("H96").Value = 4849.778
$ws.Range("I96").Value = 1736.2727
$ws.Range("J96").Value = 9742.429
$ws.Range("K96").Value = 1736.2727
$ws.Range("L96").Value = 9742.429
$ws.Range("M96").Value = -363.2727
$ws.Range("N96").Value = -12488.429

$ws.Range("H107").Value = 286.9565
$ws.Range("I107").Value = 281.92856
$ws.Range("K107").Value = 845.78568
$ws.Range("M107").Value = 1074.21432

$ws.Range("H132").Value = 20188.473
$ws.Range("I132").Value = 28180.918
$ws.Range("K132").Value = 84542.754
$ws.Range("M132").Value = -82012.754
